$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the values first
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the "header" style on B1: bold font, thin box border,
# centered horizontally, top-aligned vertically
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160

# Copy that exact formatting onto A2 so both cells share one style record
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
